$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0,0,0,0,0),
    @(0.158198511282029,2.382964774536509,0.05411120097525475,247.8342241219283,-328.1946692603154),
    @(0.1540195103665217,1.788466067858054,0.0586987627156471,303.4561756544852,-208.5430616732992),
    @(0.1376840416828386,1.205309520067421,0.05003394674851547,293.6285991707177,-257.60970680964),
    @(0.1019788141093302,1.504276112067938,0.03327306634794119,340.048837864487,-305.2886354542786),
    @(0.1167150697089643,1.685242637903506,0.03231202834206442,338.8017652129342,-300.9313775146296),
    @(0.1009261513800769,1.225863637484949,0.03162228579006612,396.8258629731951,-253.9037553476354),
    @(0.06094166116993887,0.7625901528269406,0.01849781433443191,192.9238904066669,-177.8628334406862),
    @(0.07118947121057537,1.119460267558798,0.02337542951264013,192.6881573065673,-185.8109313328811),
    @(0.0541078498096448,0.7588137205112726,0.01852106772350293,197.0010493325451,-192.5617075595327),
    @(0.03714639929402445,0.2436152250794927,0.01563910710119206,93.62070615841097,-126.5700158953849),
    @(0.02670730936268782,0.5478846551529415,0.01179336757600432,104.6479971095391,-113.8645512141335),
    @(0.0408770819726052,0.2522657329159978,0.01478029039727919,97.73063300077291,-118.9361617004882),
    @(0.01637990724578956,0.3090364102667344,0.005688285407817011,67.7750465728467,-80.43421957432557),
    @(0.01514138919410355,0.4120035099296989,0.009063658100976027,73.00548254706821,-73.37087307661103),
    @(0.01695514102691537,0.3271829186070354,0.005380600201584886,62.20673974235961,-79.38237905458644),
    @(0.01277602165241459,0.1632611229311623,0.006320965955682019,27.15632824073893,-50.78388276646678),
    @(0.009184853667081782,0.1876403846026099,0.005946669553404348,31.26249992272761,-67.37116608050415),
    @(0.002866464322252697,0.2435478807864958,0.008378269355337786,34.34349841099724,-72.29254313702833),
    @(0.002078510743094199,0.1845282344935664,0.005325596082023524,35.54178086282234,-70.26433410178562),
    @(0.01097339483756912,0.1262704788765468,0.004769659520192954,43.28600926053447,-41.65370072769352),
    @(0.006234698565036533,0.1822255452505377,0.006497059408533399,27.67473969566274,-57.40356582818455),
    @(0.00513491467248337,0.00588135787470902,0.004168230232110934,42.91778683065272,-40.41788011606975),
    @(0.005376217746727448,-0.02325644003000819,0.003758026797506845,35.20149788051975,-46.65545164026287),
    @(0.00483635662983685,0.05741663671836112,0.003730438096506507,48.08190431054824,-64.59116634311304),
    @(0.0026746514566514,0.01584425828584345,0.002063345441072744,39.84994971997929,-56.53078962697166),
    @(0.006714176797792768,-0.03158908954096089,0.004751901111793445,39.16005634655902,-51.99881467912345),
    @(-0.001693836225436336,0.2072997783691784,0.005098871097328627,41.76056848758277,-37.99410850835997),
    @(0.002330665727222756,0.08029727403436233,0.003238637763507987,38.37018805138749,-50.28915177815563),
    @(0.001851609066906718,0.06012627387253142,0.002653532385341967,32.22939598880857,-66.33103358886862),
    @(0.0100494872671646,-0.3598738674937444,0.006856646187861419,11.01814695122591,-9.608064070489965),
    @(0.0007339728299216452,-0.04715314444366168,0.001510929242360452,19.78744087512309,-59.21353538232588),
    @(0.00565511710181539,-0.1822266930381763,0.004259163515995619,20.9308915086084,-50.3774663521443),
    @(0.005019019614673753,-0.1619196163743815,0.004201990327278174,19.57587411217437,-49.4390421254724),
    @(0.001585721261110376,-0.08420676243206225,0.001781603747779381,20.80094986944983,-50.18766739303585),
    @(0.003638466906014801,-0.1331171612016337,0.002538966893398722,12.98525938162125,-63.76945084478938),
    @(0.004996495700326854,-0.1896409366880591,0.00425573512797767,16.0994107644036,-73.56742789005463),
    @(0.007419953615086644,-0.1921566973228025,0.004808294806688117,18.08096162301872,-46.98239848817071),
    @(0.01100494228349626,-0.4391384069867468,0.007817119392250431,12.56635015022701,-9.350520304598183),
    @(-0.001403690472340024,0.1212159762492618,0.002633551022753115,19.24442053694606,-49.12251450380327),
    @(0.003998759063008198,-0.07155415964064167,0.003232474042622202,19.59803963002415,-52.0593832450757)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
